$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 983
$ws.Range("F5").Value = 441
$ws.Range("F6").Value = 686
$ws.Range("F7").Value = 243
$ws.Range("F9").Value = 18
$ws.Range("F10").Value = 386
$ws.Range("F11").Value = 190
$ws.Range("F13").Value = 787
$ws.Range("F14").Value = 108
$ws.Range("F15").Value = 1942
$ws.Range("F16").Value = 444
$ws.Range("F17").Value = 6442
$ws.Range("F18").Value = 502
$ws.Range("F20").Value = 44
$ws.Range("F21").Value = 83
$ws.Range("F23").Value = 203
$ws.Range("G23").Value = 29.9
$ws.Range("F24").Value = 135

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 509
$ws.Range("F12").Value = 111
$ws.Range("F13").Value = 50
$ws.Range("F15").Value = 1
$ws.Range("F18").Value = 26

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 373
$ws.Range("F4").Value = 368

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 373
$ws.Range("F5").Value = 368
$ws.Range("F10").Value = 509
$ws.Range("F11").Value = 509
$ws.Range("F12").Value = 983
$ws.Range("F16").Value = 442
$ws.Range("F17").Value = 686
$ws.Range("F18").Value = 243
$ws.Range("F21").Value = 18
$ws.Range("F22").Value = 386
$ws.Range("F23").Value = 190
$ws.Range("F27").Value = 787
$ws.Range("F28").Value = 108
$ws.Range("F29").Value = 111
$ws.Range("F30").Value = 1942
$ws.Range("F31").Value = 444
$ws.Range("F32").Value = 6444
$ws.Range("F33").Value = 50
$ws.Range("F34").Value = 502
$ws.Range("F36").Value = 44
$ws.Range("F37").Value = 83
$ws.Range("F40").Value = 203
$ws.Range("G40").Value = 29.9
$ws.Range("F41").Value = 1
$ws.Range("F42").Value = 135
$ws.Range("F45").Value = 26
